$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Milestone Start (C28) slipped from 2019-11-18 to 2019-11-29
# Milestone End   (C29) slipped from 2019-12-02 to 2019-12-12
# (serial date numbers so the cells stay plain literals, same as the source file)
$ws.Range("C28").Value = 43798
$ws.Range("C29").Value = 43811

# Scroll the view down a bit and land the cursor on C30, matching where the
# author was working after nudging the milestone dates.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("C30").Select()
